# Excel COM-interop script: "Meters GUI section added"
#
# For each of the 5 report sheets (Activa / Reactiva Inductiva / Reactiva
# Capacitiva / Reactiva Inductiva Penalizada / Factor de Potencia) this:
#   1. Appends a new meter-reading row (account 11002006, date 29/02/2024)
#      as row 31, cloning row 30s formatting so styles/borders match.
#   2. On the sheets where day 28/02/2024 (row 30) had only partially-loaded
#      readings, fills in the finalized hourly figures.
#   3. Refreshes the sheets visible selection to cover the new used range.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Activa ---
$ws = $wb.Worksheets.Item(1)

# Clone row 30 into row 31 first so the new row inherits correct cell styles/borders
$ws.Range("A30:AA30").Copy($ws.Range("A31:AA31"))

# New row 31: meter 11002006, date 29/02/2024
$ws.Range("A31").Value = 11002006
$ws.Range("B31").Value = "29/02/2024"
$ws.Range("C31").Value = 690.59
$ws.Range("D31").Value = 596.29
$ws.Range("E31").Value = 702.08
$ws.Range("F31").Value = 768.21
$ws.Range("G31").Value = 779.71
$ws.Range("H31").Value = 933.82
$ws.Range("I31").Value = 1492.71
$ws.Range("J31").Value = 1907.86
$ws.Range("K31").Value = 2042.41
$ws.Range("L31").Value = 2068.86
$ws.Range("M31").Value = 1801.48
$ws.Range("N31").Value = 1784.24
$ws.Range("O31").Value = 1980.31
$ws.Range("P31").Value = 2036.08
$ws.Range("Q31").Value = 1872.22
$ws.Range("R31").Value = 1587.58
$ws.Range("S31").Value = 1099.98
$ws.Range("T31").Value = 1001.09
$ws.Range("U31").Value = 852.16
$ws.Range("V31").Value = 1075.26
$ws.Range("W31").Value = 1091.94
$ws.Range("X31").Value = 914.84
$ws.Range("Y31").Value = 847.55
$ws.Range("Z31").Value = 924.04
$ws.Range("AA31").Value = 30851.31

# Finalized values for row 30 (28/02/2024)
$ws.Range("D30").Value = 606.64
$ws.Range("E30").Value = 609.5
$ws.Range("F30").Value = 605.49
$ws.Range("G30").Value = 645.74
$ws.Range("H30").Value = 839.51
$ws.Range("I30").Value = 1462.81
$ws.Range("J30").Value = 1787.68
$ws.Range("K30").Value = 1932.58
$ws.Range("L30").Value = 1981.45
$ws.Range("M30").Value = 1749.17
$ws.Range("N30").Value = 1753.76
$ws.Range("O30").Value = 2015.39
$ws.Range("P30").Value = 2042.41
$ws.Range("Q30").Value = 1919.93
$ws.Range("R30").Value = 1673.84
$ws.Range("S30").Value = 1252.36
$ws.Range("T30").Value = 1096.53
$ws.Range("U30").Value = 844.68
$ws.Range("V30").Value = 1039.61
$ws.Range("W30").Value = 1083.31
$ws.Range("X30").Value = 901.6
$ws.Range("Y30").Value = 873.43
$ws.Range("Z30").Value = 841.24
$ws.Range("AA30").Value = 30262.47

# Refresh the visible selection to match the new used range
$ws.Range("A3:A31").Select() | Out-Null

# --- Sheet 2: Reactiva Inductiva ---
$ws = $wb.Worksheets.Item(2)

# Clone row 30 into row 31 first so the new row inherits correct cell styles/borders
$ws.Range("A30:AA30").Copy($ws.Range("A31:AA31"))

# New row 31: meter 11002006, date 29/02/2024
$ws.Range("A31").Value = 11002006
$ws.Range("B31").Value = "29/02/2024"
$ws.Range("C31").Value = 163.89
$ws.Range("D31").Value = 101.79
$ws.Range("E31").Value = 178.84
$ws.Range("F31").Value = 239.21
$ws.Range("G31").Value = 246.11
$ws.Range("H31").Value = 300.73
$ws.Range("I31").Value = 543.96
$ws.Range("J31").Value = 790.07
$ws.Range("K31").Value = 896.43
$ws.Range("L31").Value = 886.09
$ws.Range("M31").Value = 665.86
$ws.Range("N31").Value = 695.19
$ws.Range("O31").Value = 865.38
$ws.Range("P31").Value = 956.82
$ws.Range("Q31").Value = 1039.04
$ws.Range("R31").Value = 730.83
$ws.Range("S31").Value = 443.33
$ws.Range("T31").Value = 426.08
$ws.Range("U31").Value = 336.39
$ws.Range("V31").Value = 453.11
$ws.Range("W31").Value = 478.41
$ws.Range("X31").Value = 394.46
$ws.Range("Y31").Value = 374.91
$ws.Range("Z31").Value = 429.54
$ws.Range("AA31").Value = 12636.47

# Finalized values for row 30 (28/02/2024)
$ws.Range("D30").Value = 112.72
$ws.Range("E30").Value = 116.73
$ws.Range("F30").Value = 135.72
$ws.Range("G30").Value = 161.59
$ws.Range("H30").Value = 224.26
$ws.Range("I30").Value = 548.56
$ws.Range("J30").Value = 700.94
$ws.Range("K30").Value = 796.38
$ws.Range("L30").Value = 908.52
$ws.Range("M30").Value = 669.31
$ws.Range("N30").Value = 662.41
$ws.Range("O30").Value = 902.18
$ws.Range("P30").Value = 927.48
$ws.Range("Q30").Value = 779.14
$ws.Range("R30").Value = 608.36
$ws.Range("S30").Value = 443.92
$ws.Range("T30").Value = 392.17
$ws.Range("U30").Value = 285.21
$ws.Range("V30").Value = 389.86
$ws.Range("W30").Value = 413.43
$ws.Range("X30").Value = 299.01
$ws.Range("Y30").Value = 295.56
$ws.Range("Z30").Value = 275.44
$ws.Range("AA30").Value = 11253.61

# Refresh the visible selection to match the new used range
$ws.Range("A3:A31").Select() | Out-Null

# --- Sheet 3: Reactiva Capacitiva ---
$ws = $wb.Worksheets.Item(3)

# Clone row 30 into row 31 first so the new row inherits correct cell styles/borders
$ws.Range("A30:AA30").Copy($ws.Range("A31:AA31"))

# New row 31: meter 11002006, date 29/02/2024
$ws.Range("A31").Value = 11002006
$ws.Range("B31").Value = "29/02/2024"
$ws.Range("C31").Value = 0.0
$ws.Range("D31").Value = 0.0
$ws.Range("E31").Value = 0.0
$ws.Range("F31").Value = 0.0
$ws.Range("G31").Value = 0.0
$ws.Range("H31").Value = 0.0
$ws.Range("I31").Value = 0.0
$ws.Range("J31").Value = 0.0
$ws.Range("K31").Value = 0.0
$ws.Range("L31").Value = 0.0
$ws.Range("M31").Value = 0.0
$ws.Range("N31").Value = 0.0
$ws.Range("O31").Value = 0.0
$ws.Range("P31").Value = 0.0
$ws.Range("Q31").Value = 0.0
$ws.Range("R31").Value = 0.0
$ws.Range("S31").Value = 0.0
$ws.Range("T31").Value = 0.0
$ws.Range("U31").Value = 0.0
$ws.Range("V31").Value = 0.0
$ws.Range("W31").Value = 0.0
$ws.Range("X31").Value = 0.0
$ws.Range("Y31").Value = 0.0
$ws.Range("Z31").Value = 0.0
$ws.Range("AA31").Value = 0.0

# Refresh the visible selection to match the new used range
$ws.Range("A3:A31").Select() | Out-Null

# --- Sheet 4: Reactiva Inductiva Penalizada ---
$ws = $wb.Worksheets.Item(4)

# Clone row 30 into row 31 first so the new row inherits correct cell styles/borders
$ws.Range("A30:AA30").Copy($ws.Range("A31:AA31"))

# New row 31: meter 11002006, date 29/02/2024
$ws.Range("A31").Value = 11002006
$ws.Range("B31").Value = "29/02/2024"
$ws.Range("C31").Value = 0.0
$ws.Range("D31").Value = 0.0
$ws.Range("E31").Value = 0.0
$ws.Range("F31").Value = 0.0
$ws.Range("G31").Value = 0.0
$ws.Range("H31").Value = 0.0
$ws.Range("I31").Value = 0.0
$ws.Range("J31").Value = 0.0
$ws.Range("K31").Value = 0.0
$ws.Range("L31").Value = 0.0
$ws.Range("M31").Value = 0.0
$ws.Range("N31").Value = 0.0
$ws.Range("O31").Value = 0.0
$ws.Range("P31").Value = 0.0
$ws.Range("Q31").Value = 102.93
$ws.Range("R31").Value = 0.0
$ws.Range("S31").Value = 0.0
$ws.Range("T31").Value = 0.0
$ws.Range("U31").Value = 0.0
$ws.Range("V31").Value = 0.0
$ws.Range("W31").Value = 0.0
$ws.Range("X31").Value = 0.0
$ws.Range("Y31").Value = 0.0
$ws.Range("Z31").Value = 0.0
$ws.Range("AA31").Value = 102.93

# Refresh the visible selection to match the new used range
$ws.Range("A3:A31").Select() | Out-Null

# --- Sheet 5: Factor de Potencia ---
$ws = $wb.Worksheets.Item(5)

# Clone row 30 into row 31 first so the new row inherits correct cell styles/borders
$ws.Range("A30:AA30").Copy($ws.Range("A31:AA31"))

# New row 31: meter 11002006, date 29/02/2024
$ws.Range("A31").Value = 11002006
$ws.Range("B31").Value = "29/02/2024"
$ws.Range("C31").Value = 0.97
$ws.Range("D31").Value = 0.99
$ws.Range("E31").Value = 0.97
$ws.Range("F31").Value = 0.95
$ws.Range("G31").Value = 0.95
$ws.Range("H31").Value = 0.95
$ws.Range("I31").Value = 0.94
$ws.Range("J31").Value = 0.92
$ws.Range("K31").Value = 0.92
$ws.Range("L31").Value = 0.92
$ws.Range("M31").Value = 0.94
$ws.Range("N31").Value = 0.93
$ws.Range("O31").Value = 0.92
$ws.Range("P31").Value = 0.91
$ws.Range("Q31").Value = 0.87
$ws.Range("R31").Value = 0.91
$ws.Range("S31").Value = 0.93
$ws.Range("T31").Value = 0.92
$ws.Range("U31").Value = 0.93
$ws.Range("V31").Value = 0.92
$ws.Range("W31").Value = 0.92
$ws.Range("X31").Value = 0.91
$ws.Range("Y31").Value = 0.91
$ws.Range("Z31").Value = 0.91
$ws.Range("AA31").Value = 0.93

# Finalized values for row 30 (28/02/2024)
$ws.Range("E30").Value = 0.98
$ws.Range("F30").Value = 0.98
$ws.Range("G30").Value = 0.97
$ws.Range("H30").Value = 0.97
$ws.Range("I30").Value = 0.94
$ws.Range("J30").Value = 0.93
$ws.Range("K30").Value = 0.92
$ws.Range("L30").Value = 0.91
$ws.Range("M30").Value = 0.93
$ws.Range("N30").Value = 0.94
$ws.Range("O30").Value = 0.91
$ws.Range("P30").Value = 0.91
$ws.Range("Q30").Value = 0.93
$ws.Range("R30").Value = 0.94
$ws.Range("S30").Value = 0.94
$ws.Range("T30").Value = 0.94
$ws.Range("U30").Value = 0.95
$ws.Range("V30").Value = 0.94
$ws.Range("W30").Value = 0.93
$ws.Range("X30").Value = 0.95
$ws.Range("Y30").Value = 0.95
$ws.Range("Z30").Value = 0.95
$ws.Range("AA30").Value = 0.94

# Refresh the visible selection to match the new used range
$ws.Range("A3:A31").Select() | Out-Null
